$d = $word.ActiveDocument

# Locate the end of the sentence that currently precedes the page break
# ("... to manage the damage outcome.") so we can split the paragraph
# right after it.
$target = $d.Content
$found = $target.Find.Execute("to manage the damage outcome.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor sentence for the edit."
}

# Remember which paragraph holds the sentence before we mutate anything.
$paraIndex = $target.Paragraphs.First.Index

# Collapse to the end of the match and split the paragraph there. This
# leaves the trailing page-break run alone in a brand new paragraph,
# exactly like pressing Enter right after "...damage outcome." and
# before the page break.
$target.Collapse(0)
$target.InsertParagraphAfter()

# The paragraph that used to hold "...damage outcome.<page break>" is now
# two paragraphs: $paraIndex (the sentence) and $paraIndex + 1 (just the
# page-break run).
$newPara = $d.Paragraphs.Item($paraIndex + 1)

# Insert the new sentence text at the very start of that paragraph, ahead
# of the page-break run.
$insertion = $newPara.Range
$insertion.Collapse(1)
$insertion.InsertBefore("All of these scripts can be found in the ComboSystem namespace.")

# Split the trailing period into its own run, matching the document's
# existing convention elsewhere (every other sentence in this file ends
# with its period as a separate run). Toggling a direct character
# property on just that one character and then clearing it again forces
# the run boundary without altering the visible formatting.
$paraRange = $d.Paragraphs.Item($paraIndex + 1).Range
$text = $paraRange.Text
$periodIndex = $text.Length - 3  # char right before the page-break + paragraph mark
$periodPos = $paraRange.Start + $periodIndex
$periodRange = $d.Range($periodPos, $periodPos + 1)
$periodRange.Bold = $true
$periodRange.Bold = $false
